$wb = $excel.ActiveWorkbook

# --- Typography sheet: fill in the "Widget Wildcard Characters" example for the Default font row ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("G4").Value = "0123456789-"

# --- Translation sheet: add birthday-collision example rows ---
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("B4").Value = "SingleUseId3"
$ws2.Range("C4").Value = "Default"
$ws2.Range("D4").Value = "Center"
$ws2.Range("E4").Value = "LTR"
$ws2.Range("F4").Value = "<value> at <value>"

$ws2.Range("B5").Value = "SingleUseId4"
$ws2.Range("C5").Value = "Default"
$ws2.Range("D5").Value = "Left"
$ws2.Range("E5").Value = "LTR"
$ws2.Range("F5").NumberFormat = "@"
$ws2.Range("F5").Value = "0"
$ws2.Range("F5").Style = "Normal"

$ws2.Range("B6").Value = "SingleUseId5"
$ws2.Range("C6").Value = "Default"
$ws2.Range("D6").Value = "Left"
$ws2.Range("E6").Value = "LTR"
$ws2.Range("F6").NumberFormat = "@"
$ws2.Range("F6").Value = "0"
$ws2.Range("F6").Style = "Normal"
